$d = $word.ActiveDocument

# Locate the run of text that needs updating.
$rng = $d.Content
$found = $rng.Find.Execute("Marketing Science Institute (MSI) Scholar, 2023", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not find 'Marketing Science Institute (MSI) Scholar, 2023' in the document"
}

# $rng now spans the found text ("Marketing Science Institute (MSI) Scholar, 2023").
# We need to change the trailing "3" to "2" (2023 -> 2022), and leave that final
# character in its own run (split off from the rest of the sentence), matching how
# Word splits a run when the tail of it is edited/retyped.
$lastChar = $d.Range($rng.End - 1, $rng.End)

# Force a run split at this boundary without altering the visible text: toggling a
# character-formatting property on/off is a no-op visually but causes Word's run
# model to carve this single character out into its own <w:r>.
$lastChar.Font.Bold = 1
$lastChar.Font.Bold = 0

# Re-acquire the now-isolated trailing run and give it distinct formatting momentarily
# so replacing its text doesn't get re-coalesced with the preceding run, then restore
# the formatting to match (Bold off, same as the rest of the sentence).
$lastChar.Font.Bold = 1
$lastChar.Text = "2"
$lastChar.Font.Bold = 0
